# Applies the scheduled market-data refresh: updates current/leve price and
# profit figures across the ALC/BSM/CRP/CUL/GSM/LTW/WVR sheets, and removes
# now-empty profit cells where the refreshed data no longer produces a value.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1216.2354
$ws.Range("I33").Value = 245.06667
$ws.Range("K33").Value = 245.06667
$ws.Range("M33").Value = -16.06666999999999
$ws.Range("H49").Value = 430
$ws.Range("I49").Value = 300
$ws.Range("K49").Value = 900
$ws.Range("M49").Value = -764
$ws.Range("H62").Value = 2098.5386
$ws.Range("I62").Value = 1763.7273
$ws.Range("J62").Value = 3940
$ws.Range("K62").Value = 1763.7273
$ws.Range("L62").Value = 3940
$ws.Range("M62").Value = -1139.7273
$ws.Range("N62").Value = -5188
$ws.Range("H64").Value = 34187.344
$ws.Range("I64").Value = 79631.46000000001
$ws.Range("J64").Value = 3094
$ws.Range("K64").Value = 79631.46000000001
$ws.Range("L64").Value = 3094
$ws.Range("M64").Value = -79383.46000000001
$ws.Range("N64").Value = -3590
$ws.Range("H65").Value = 2098.5386
$ws.Range("I65").Value = 1763.7273
$ws.Range("J65").Value = 3940
$ws.Range("K65").Value = 8818.636500000001
$ws.Range("L65").Value = 19700
$ws.Range("M65").Value = -5698.636500000001
$ws.Range("N65").Value = -25940
$ws.Range("H67").Value = 34187.344
$ws.Range("I67").Value = 79631.46000000001
$ws.Range("J67").Value = 3094
$ws.Range("K67").Value = 79631.46000000001
$ws.Range("L67").Value = 3094
$ws.Range("M67").Value = -78773.46000000001
$ws.Range("N67").Value = -4810
$ws.Range("H111").Value = 6259238
$ws.Range("J111").Value = 25002498
$ws.Range("L111").Value = 75007494
$ws.Range("N111").Value = -75013628
$ws.Range("H112").Value = 1081.4375
$ws.Range("J112").Value = 1165.2142
$ws.Range("L112").Value = 3495.6426
$ws.Range("N112").Value = -5711.642599999999
$ws.Range("H113").Value = 2412
$ws.Range("I113").Value = 2824.5
$ws.Range("J113").Value = 1999.5
$ws.Range("K113").Value = 2824.5
$ws.Range("L113").Value = 1999.5
$ws.Range("M113").Value = 429.5
$ws.Range("N113").Value = -8507.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2472.5818
$ws.Range("I31").Value = 1640.4762
$ws.Range("J31").Value = 2986.5293
$ws.Range("K31").Value = 1640.4762
$ws.Range("L31").Value = 2986.5293
$ws.Range("M31").Value = -1345.4762
$ws.Range("N31").Value = -3576.5293
$ws.Range("H34").Value = 2472.5818
$ws.Range("I34").Value = 1640.4762
$ws.Range("J34").Value = 2986.5293
$ws.Range("K34").Value = 1640.4762
$ws.Range("L34").Value = 2986.5293
$ws.Range("M34").Value = -1438.4762
$ws.Range("N34").Value = -3390.5293
$ws.Range("H58").Value = 2439.0476
$ws.Range("I58").Value = 1622.9
$ws.Range("K58").Value = 1622.9
$ws.Range("M58").Value = -1419.9
$ws.Range("H132").Value = 6213.45
$ws.Range("I132").Value = 6721.4165
$ws.Range("K132").Value = 20164.2495
$ws.Range("M132").Value = -17634.2495
$ws.Range("H136").Value = 2439.0476
$ws.Range("I136").Value = 1622.9
$ws.Range("K136").Value = 4868.700000000001
$ws.Range("M136").Value = -2318.700000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1153.3673
$ws.Range("I5").Value = 1026.2941
$ws.Range("J5").Value = 1220.875
$ws.Range("K5").Value = 3078.8823
$ws.Range("L5").Value = 3662.625
$ws.Range("M5").Value = -2966.8823
$ws.Range("N5").Value = -3886.625
$ws.Range("H120").Value = 8900
$ws.Range("I120").Value = 8900
$ws.Range("K120").Value = 26700
$ws.Range("M120").Value = -21862
$ws.Range("H122").Value = 573.7778
$ws.Range("I122").Value = 549.8333
$ws.Range("K122").Value = 4948.4997
$ws.Range("M122").Value = -2498.4997
$ws.Range("H131").Value = 751.9899
$ws.Range("I131").Value = 289.9
$ws.Range("J131").Value = 803.9101000000001
$ws.Range("K131").Value = 869.6999999999999
$ws.Range("L131").Value = 2411.7303
$ws.Range("M131").Value = 4170.3
$ws.Range("N131").Value = -12491.7303
$ws.Range("H135").Value = 1153.3673
$ws.Range("I135").Value = 1026.2941
$ws.Range("J135").Value = 1220.875
$ws.Range("K135").Value = 9236.6469
$ws.Range("L135").Value = 10987.875
$ws.Range("M135").Value = -6701.6469
$ws.Range("N135").Value = -16057.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1414.5
$ws.Range("I113").Value = 970.1429000000001
$ws.Range("J113").Value = 1858.8572
$ws.Range("K113").Value = 970.1429000000001
$ws.Range("L113").Value = 1858.8572
$ws.Range("M113").Value = 1199.8571
$ws.Range("N113").Value = -6198.8572

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 64781.562
$ws.Range("I40").Value = 144685.72
$ws.Range("J40").Value = 2633.889
$ws.Range("K40").Value = 144685.72
$ws.Range("L40").Value = 2633.889
$ws.Range("M40").Value = -144549.72
$ws.Range("N40").Value = -2905.889
$ws.Range("H46").Value = 1265774.6
$ws.Range("I46").Value = 545
$ws.Range("J46").Value = 1687517.9
$ws.Range("K46").Value = 545
$ws.Range("L46").Value = 1687517.9
$ws.Range("M46").Value = -357
$ws.Range("N46").Value = -1687893.9
$ws.Range("H61").Value = 1602.6154
$ws.Range("I61").Value = 1550.5
$ws.Range("J61").Value = 1686
$ws.Range("K61").Value = 1550.5
$ws.Range("L61").Value = 1686
$ws.Range("M61").Value = -1348.5
$ws.Range("N61").Value = -2090
$ws.Range("H113").Value = 1602.6154
$ws.Range("I113").Value = 1550.5
$ws.Range("J113").Value = 1686
$ws.Range("K113").Value = 1550.5
$ws.Range("L113").Value = 1686
$ws.Range("M113").Value = 619.5
$ws.Range("N113").Value = -6026

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 286155.56
$ws.Range("I107").Value = 630
$ws.Range("J107").Value = 500299.75
$ws.Range("K107").Value = 1890
$ws.Range("L107").Value = 1500899.25
$ws.Range("M107").Value = 30
$ws.Range("N107").Value = -1504739.25
$ws.Range("H113").Value = 997
$ws.Range("I113").Value = 956.125
$ws.Range("J113").Value = 1043.7142
$ws.Range("K113").Value = 2868.375
$ws.Range("L113").Value = 3131.1426
$ws.Range("M113").Value = -698.375
$ws.Range("N113").Value = -7471.142599999999
$ws.Range("H136").Value = 1124.64
$ws.Range("I136").Value = 383.58536
$ws.Range("J136").Value = 4500.5557
$ws.Range("K136").Value = 1150.75608
$ws.Range("L136").Value = 13501.6671
$ws.Range("M136").Value = 1399.24392
$ws.Range("N136").Value = -18601.6671
